$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant data: rows 2-37, columns A-T

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Dag1"
$ws.Range("C2").Value2 = "Lama4"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 21.0077345
$ws.Range("H2").Value2 = 42.015469
$ws.Range("I2").Value2 = 0.07974258627637139
$ws.Range("J2").Value2 = 0.07240574124537677
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 104.0097425
$ws.Range("N2").Value2 = 208.019485
$ws.Range("O2").Value2 = 0.4223541729035338
$ws.Range("P2").Value2 = 0.3348767214268842
$ws.Range("Q2").Value2 = 2185.009055853366
$ws.Range("R2").Value2 = 8740.036223413465
$ws.Range("S2").Value2 = 0.03367961407194552
$ws.Range("T2").Value2 = 0.02424699724073509

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Dag1"
$ws.Range("C3").Value2 = "Lama4"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 21.0077345
$ws.Range("H3").Value2 = 42.015469
$ws.Range("I3").Value2 = 0.07974258627637139
$ws.Range("J3").Value2 = 0.07240574124537677
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 58.20636966666666
$ws.Range("N3").Value2 = 174.619109
$ws.Range("O3").Value2 = 0.2363596190835905
$ws.Range("P3").Value2 = 0.2811076795060988
$ws.Range("Q3").Value2 = 1222.783960166186
$ws.Range("R3").Value2 = 7336.703760997119
$ws.Range("S3").Value2 = 0.01884792731702349
$ws.Range("T3").Value2 = 0.02035380990440689

# Row 4
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Dag1"
$ws.Range("C4").Value2 = "Lama4"
$ws.Range("D4").Value2 = "M1"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 21.0077345
$ws.Range("H4").Value2 = 42.015469
$ws.Range("I4").Value2 = 0.07974258627637139
$ws.Range("J4").Value2 = 0.07240574124537677
$ws.Range("K4").Value2 = 2
$ws.Range("L4").Value2 = 0.6666666666666666
$ws.Range("M4").Value2 = 0.03002933333333333
$ws.Range("N4").Value2 = 0.090088
$ws.Range("O4").Value2 = 0.0001219406368864389
$ws.Range("P4").Value2 = 0.0001450266742074914
$ws.Range("Q4").Value2 = 0.6308482618786666
$ws.Range("R4").Value2 = 3.785089571272
$ws.Range("S4").Value2 = 0.000009723861757512532
$ws.Range("T4").Value2 = 0.00001050076384634518

# Row 5
$ws.Range("A5").Value2 = "ECs"
$ws.Range("B5").Value2 = "Dag1"
$ws.Range("C5").Value2 = "Lama4"
$ws.Range("D5").Value2 = "M2"
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 21.0077345
$ws.Range("H5").Value2 = 42.015469
$ws.Range("I5").Value2 = 0.07974258627637139
$ws.Range("J5").Value2 = 0.07240574124537677
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.1136856666666667
$ws.Range("N5").Value2 = 0.341057
$ws.Range("O5").Value2 = 0.0004616453666923254
$ws.Range("P5").Value2 = 0.0005490449607626364
$ws.Range("Q5").Value2 = 2.388278301788833
$ws.Range("R5").Value2 = 14.329669810733
$ws.Range("S5").Value2 = 0.00003681279548254986
$ws.Range("T5").Value2 = 0.00003975400736105749

# Row 6
$ws.Range("A6").Value2 = "ECs"
$ws.Range("B6").Value2 = "Dag1"
$ws.Range("C6").Value2 = "Lama4"
$ws.Range("D6").Value2 = "Neutro"
$ws.Range("E6").Value2 = 2
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 21.0077345
$ws.Range("H6").Value2 = 42.015469
$ws.Range("I6").Value2 = 0.07974258627637139
$ws.Range("J6").Value2 = 0.07240574124537677
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 70.30838766666666
$ws.Range("N6").Value2 = 210.925163
$ws.Range("O6").Value2 = 0.285502494356584
$ws.Range("P6").Value2 = 0.3395543790134426
$ws.Range("Q6").Value2 = 1477.019941224408
$ws.Range("R6").Value2 = 8862.119647346446
$ws.Range("S6").Value2 = 0.02276670728834914
$ws.Range("T6").Value2 = 0.02458568650558192

# Row 7
$ws.Range("A7").Value2 = "ECs"
$ws.Range("B7").Value2 = "Dag1"
$ws.Range("C7").Value2 = "Lama4"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 2
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 21.0077345
$ws.Range("H7").Value2 = 42.015469
$ws.Range("I7").Value2 = 0.07974258627637139
$ws.Range("J7").Value2 = 0.07240574124537677
$ws.Range("K7").Value2 = 2
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 13.5936885
$ws.Range("N7").Value2 = 27.187377
$ws.Range("O7").Value2 = 0.05520012765271271
$ws.Range("P7").Value2 = 0.04376714841860452
$ws.Range("Q7").Value2 = 285.5725988837032
$ws.Range("R7").Value2 = 1142.290395534813
$ws.Range("S7").Value2 = 0.004401800941813158
$ws.Range("T7").Value2 = 0.00316899282344548

# Row 8
$ws.Range("A8").Value2 = "FAPs"
$ws.Range("B8").Value2 = "Dag1"
$ws.Range("C8").Value2 = "Lama4"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 23.794572
$ws.Range("H8").Value2 = 71.38371599999999
$ws.Range("I8").Value2 = 0.09032105344911566
$ws.Range("J8").Value2 = 0.1230163792728212
$ws.Range("K8").Value2 = 2
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 104.0097425
$ws.Range("N8").Value2 = 208.019485
$ws.Range("O8").Value2 = 0.4223541729035338
$ws.Range("P8").Value2 = 0.3348767214268842
$ws.Range("Q8").Value2 = 2474.86730661771
$ws.Range("R8").Value2 = 14849.20383970626
$ws.Range("S8").Value2 = 0.03814747382527711
$ws.Range("T8").Value2 = 0.04119532177268847

# Row 9
$ws.Range("A9").Value2 = "FAPs"
$ws.Range("B9").Value2 = "Dag1"
$ws.Range("C9").Value2 = "Lama4"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 23.794572
$ws.Range("H9").Value2 = 71.38371599999999
$ws.Range("I9").Value2 = 0.09032105344911566
$ws.Range("J9").Value2 = 0.1230163792728212
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 58.20636966666666
$ws.Range("N9").Value2 = 174.619109
$ws.Range("O9").Value2 = 0.2363596190835905
$ws.Range("P9").Value2 = 0.2811076795060988
$ws.Range("Q9").Value2 = 1384.995653892116
$ws.Range("R9").Value2 = 12464.96088502904
$ws.Range("S9").Value2 = 0.02134824978846159
$ws.Range("T9").Value2 = 0.03458084891862491

# Row 10
$ws.Range("A10").Value2 = "FAPs"
$ws.Range("B10").Value2 = "Dag1"
$ws.Range("C10").Value2 = "Lama4"
$ws.Range("D10").Value2 = "M1"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 23.794572
$ws.Range("H10").Value2 = 71.38371599999999
$ws.Range("I10").Value2 = 0.09032105344911566
$ws.Range("J10").Value2 = 0.1230163792728212
$ws.Range("K10").Value2 = 2
$ws.Range("L10").Value2 = 0.6666666666666666
$ws.Range("M10").Value2 = 0.03002933333333333
$ws.Range("N10").Value2 = 0.090088
$ws.Range("O10").Value2 = 0.0001219406368864389
$ws.Range("P10").Value2 = 0.0001450266742074914
$ws.Range("Q10").Value2 = 0.714535134112
$ws.Range("R10").Value2 = 6.430816207007999
$ws.Range("S10").Value2 = 0.00001101380678183926
$ws.Range("T10").Value2 = 0.00001784065635898463

# Row 11
$ws.Range("A11").Value2 = "FAPs"
$ws.Range("B11").Value2 = "Dag1"
$ws.Range("C11").Value2 = "Lama4"
$ws.Range("D11").Value2 = "M2"
$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 23.794572
$ws.Range("H11").Value2 = 71.38371599999999
$ws.Range("I11").Value2 = 0.09032105344911566
$ws.Range("J11").Value2 = 0.1230163792728212
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 0.1136856666666667
$ws.Range("N11").Value2 = 0.341057
$ws.Range("O11").Value2 = 0.0004616453666923254
$ws.Range("P11").Value2 = 0.0005490449607626364
$ws.Range("Q11").Value2 = 2.705101780868
$ws.Range("R11").Value2 = 24.345916027812
$ws.Range("S11").Value2 = 0.00004169629583955412
$ws.Range("T11").Value2 = 0.0000675415231310077

# Row 12
$ws.Range("A12").Value2 = "FAPs"
$ws.Range("B12").Value2 = "Dag1"
$ws.Range("C12").Value2 = "Lama4"
$ws.Range("D12").Value2 = "Neutro"
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 23.794572
$ws.Range("H12").Value2 = 71.38371599999999
$ws.Range("I12").Value2 = 0.09032105344911566
$ws.Range("J12").Value2 = 0.1230163792728212
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 70.30838766666666
$ws.Range("N12").Value2 = 210.925163
$ws.Range("O12").Value2 = 0.285502494356584
$ws.Range("P12").Value2 = 0.3395543790134426
$ws.Range("Q12").Value2 = 1672.957992538412
$ws.Range("R12").Value2 = 15056.62193284571
$ws.Range("S12").Value2 = 0.02578688605263687
$ws.Range("T12").Value2 = 0.04177075027246493

# Row 13
$ws.Range("A13").Value2 = "FAPs"
$ws.Range("B13").Value2 = "Dag1"
$ws.Range("C13").Value2 = "Lama4"
$ws.Range("D13").Value2 = "sCs"
$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 23.794572
$ws.Range("H13").Value2 = 71.38371599999999
$ws.Range("I13").Value2 = 0.09032105344911566
$ws.Range("J13").Value2 = 0.1230163792728212
$ws.Range("K13").Value2 = 2
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 13.5936885
$ws.Range("N13").Value2 = 27.187377
$ws.Range("O13").Value2 = 0.05520012765271271
$ws.Range("P13").Value2 = 0.04376714841860452
$ws.Range("Q13").Value2 = 323.455999758822
$ws.Range("R13").Value2 = 1940.735998552932
$ws.Range("S13").Value2 = 0.004985733680118672
$ws.Range("T13").Value2 = 0.005384076129552911

# Row 14
$ws.Range("A14").Value2 = "M1"
$ws.Range("B14").Value2 = "Dag1"
$ws.Range("C14").Value2 = "Lama4"
$ws.Range("D14").Value2 = "ECs"
$ws.Range("E14").Value2 = 3
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 7.255376666666667
$ws.Range("H14").Value2 = 21.76613
$ws.Range("I14").Value2 = 0.02754045181831666
$ws.Range("J14").Value2 = 0.03750982231551986
$ws.Range("K14").Value2 = 2
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 104.0097425
$ws.Range("N14").Value2 = 208.019485
$ws.Range("O14").Value2 = 0.4223541729035338
$ws.Range("P14").Value2 = 0.3348767214268842
$ws.Range("Q14").Value2 = 754.6298588405084
$ws.Range("R14").Value2 = 4527.77915304305
$ws.Range("S14").Value2 = 0.01163182474911476
$ws.Range("T14").Value2 = 0.01256116631832627

# Row 15
$ws.Range("A15").Value2 = "M1"
$ws.Range("B15").Value2 = "Dag1"
$ws.Range("C15").Value2 = "Lama4"
$ws.Range("D15").Value2 = "FAPs"
$ws.Range("E15").Value2 = 3
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 7.255376666666667
$ws.Range("H15").Value2 = 21.76613
$ws.Range("I15").Value2 = 0.02754045181831666
$ws.Range("J15").Value2 = 0.03750982231551986
$ws.Range("K15").Value2 = 3
$ws.Range("L15").Value2 = 1
$ws.Range("M15").Value2 = 58.20636966666666
$ws.Range("N15").Value2 = 174.619109
$ws.Range("O15").Value2 = 0.2363596190835905
$ws.Range("P15").Value2 = 0.2811076795060988
$ws.Range("Q15").Value2 = 422.3091363309077
$ws.Range("R15").Value2 = 3800.78222697817
$ws.Range("S15").Value2 = 0.006509450701167302
$ws.Range("T15").Value2 = 0.01054429910980187

# Row 16
$ws.Range("A16").Value2 = "M1"
$ws.Range("B16").Value2 = "Dag1"
$ws.Range("C16").Value2 = "Lama4"
$ws.Range("D16").Value2 = "M1"
$ws.Range("E16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 7.255376666666667
$ws.Range("H16").Value2 = 21.76613
$ws.Range("I16").Value2 = 0.02754045181831666
$ws.Range("J16").Value2 = 0.03750982231551986
$ws.Range("K16").Value2 = 2
$ws.Range("L16").Value2 = 0.6666666666666666
$ws.Range("M16").Value2 = 0.03002933333333333
$ws.Range("N16").Value2 = 0.090088
$ws.Range("O16").Value2 = 0.0001219406368864389
$ws.Range("P16").Value2 = 0.0001450266742074914
$ws.Range("Q16").Value2 = 0.2178741243822222
$ws.Range("R16").Value2 = 1.96086711944
$ws.Range("S16").Value2 = 0.000003358300234865819
$ws.Range("T16").Value2 = 0.000005439924780533788

# Row 17
$ws.Range("A17").Value2 = "M1"
$ws.Range("B17").Value2 = "Dag1"
$ws.Range("C17").Value2 = "Lama4"
$ws.Range("D17").Value2 = "M2"
$ws.Range("E17").Value2 = 3
$ws.Range("F17").Value2 = 1
$ws.Range("G17").Value2 = 7.255376666666667
$ws.Range("H17").Value2 = 21.76613
$ws.Range("I17").Value2 = 0.02754045181831666
$ws.Range("J17").Value2 = 0.03750982231551986
$ws.Range("K17").Value2 = 3
$ws.Range("L17").Value2 = 1
$ws.Range("M17").Value2 = 0.1136856666666667
$ws.Range("N17").Value2 = 0.341057
$ws.Range("O17").Value2 = 0.0004616453666923254
$ws.Range("P17").Value2 = 0.0005490449607626364
$ws.Range("Q17").Value2 = 0.8248323332677778
$ws.Range("R17").Value2 = 7.42349099941
$ws.Range("S17").Value2 = 0.00001271392197853911
$ws.Range("T17").Value2 = 0.00002059457892143806

# Row 18
$ws.Range("A18").Value2 = "M1"
$ws.Range("B18").Value2 = "Dag1"
$ws.Range("C18").Value2 = "Lama4"
$ws.Range("D18").Value2 = "Neutro"
$ws.Range("E18").Value2 = 3
$ws.Range("F18").Value2 = 1
$ws.Range("G18").Value2 = 7.255376666666667
$ws.Range("H18").Value2 = 21.76613
$ws.Range("I18").Value2 = 0.02754045181831666
$ws.Range("J18").Value2 = 0.03750982231551986
$ws.Range("K18").Value2 = 3
$ws.Range("L18").Value2 = 1
$ws.Range("M18").Value2 = 70.30838766666666
$ws.Range("N18").Value2 = 210.925163
$ws.Range("O18").Value2 = 0.285502494356584
$ws.Range("P18").Value2 = 0.3395543790134426
$ws.Range("Q18").Value2 = 510.1138353476878
$ws.Range("R18").Value2 = 4591.02451812919
$ws.Range("S18").Value2 = 0.007862867689836727
$ws.Range("T18").Value2 = 0.01273662442325092

# Row 19
$ws.Range("A19").Value2 = "M1"
$ws.Range("B19").Value2 = "Dag1"
$ws.Range("C19").Value2 = "Lama4"
$ws.Range("D19").Value2 = "sCs"
$ws.Range("E19").Value2 = 3
$ws.Range("F19").Value2 = 1
$ws.Range("G19").Value2 = 7.255376666666667
$ws.Range("H19").Value2 = 21.76613
$ws.Range("I19").Value2 = 0.02754045181831666
$ws.Range("J19").Value2 = 0.03750982231551986
$ws.Range("K19").Value2 = 2
$ws.Range("L19").Value2 = 1
$ws.Range("M19").Value2 = 13.5936885
$ws.Range("N19").Value2 = 27.187377
$ws.Range("O19").Value2 = 0.05520012765271271
$ws.Range("P19").Value2 = 0.04376714841860452
$ws.Range("Q19").Value2 = 98.627330356835
$ws.Range("R19").Value2 = 591.76398214101
$ws.Range("S19").Value2 = 0.001520236455984463
$ws.Range("T19").Value2 = 0.001641697960438842

# Row 20
$ws.Range("A20").Value2 = "M2"
$ws.Range("B20").Value2 = "Dag1"
$ws.Range("C20").Value2 = "Lama4"
$ws.Range("D20").Value2 = "ECs"
$ws.Range("E20").Value2 = 3
$ws.Range("F20").Value2 = 1
$ws.Range("G20").Value2 = 6.849529333333333
$ws.Range("H20").Value2 = 20.548588
$ws.Range("I20").Value2 = 0.02599990892953593
$ws.Range("J20").Value2 = 0.03541161817534047
$ws.Range("K20").Value2 = 2
$ws.Range("L20").Value2 = 1
$ws.Range("M20").Value2 = 104.0097425
$ws.Range("N20").Value2 = 208.019485
$ws.Range("O20").Value2 = 0.4223541729035338
$ws.Range("P20").Value2 = 0.3348767214268842
$ws.Range("Q20").Value2 = 712.4177822061966
$ws.Range("R20").Value2 = 4274.50669323718
$ws.Range("S20").Value2 = 0.01098117003150135
$ws.Range("T20").Value2 = 0.01185852659497868

# Row 21
$ws.Range("A21").Value2 = "M2"
$ws.Range("B21").Value2 = "Dag1"
$ws.Range("C21").Value2 = "Lama4"
$ws.Range("D21").Value2 = "FAPs"
$ws.Range("E21").Value2 = 3
$ws.Range("F21").Value2 = 1
$ws.Range("G21").Value2 = 6.849529333333333
$ws.Range("H21").Value2 = 20.548588
$ws.Range("I21").Value2 = 0.02599990892953593
$ws.Range("J21").Value2 = 0.03541161817534047
$ws.Range("K21").Value2 = 3
$ws.Range("L21").Value2 = 1
$ws.Range("M21").Value2 = 58.20636966666666
$ws.Range("N21").Value2 = 174.619109
$ws.Range("O21").Value2 = 0.2363596190835905
$ws.Range("P21").Value2 = 0.2811076795060988
$ws.Range("Q21").Value2 = 398.6862364186768
$ws.Range("R21").Value2 = 3588.176127768091
$ws.Range("S21").Value2 = 0.006145328570793155
$ws.Range("T21").Value2 = 0.00995447781282595

# Row 22
$ws.Range("A22").Value2 = "M2"
$ws.Range("B22").Value2 = "Dag1"
$ws.Range("C22").Value2 = "Lama4"
$ws.Range("D22").Value2 = "M1"
$ws.Range("E22").Value2 = 3
$ws.Range("F22").Value2 = 1
$ws.Range("G22").Value2 = 6.849529333333333
$ws.Range("H22").Value2 = 20.548588
$ws.Range("I22").Value2 = 0.02599990892953593
$ws.Range("J22").Value2 = 0.03541161817534047
$ws.Range("K22").Value2 = 2
$ws.Range("L22").Value2 = 0.6666666666666666
$ws.Range("M22").Value2 = 0.03002933333333333
$ws.Range("N22").Value2 = 0.090088
$ws.Range("O22").Value2 = 0.0001219406368864389
$ws.Range("P22").Value2 = 0.0001450266742074914
$ws.Range("Q22").Value2 = 0.2056867995271111
$ws.Range("R22").Value2 = 1.851181195744
$ws.Range("S22").Value2 = 0.000003170445453857022
$ws.Range("T22").Value2 = 0.000005135629212275182

# Row 23
$ws.Range("A23").Value2 = "M2"
$ws.Range("B23").Value2 = "Dag1"
$ws.Range("C23").Value2 = "Lama4"
$ws.Range("D23").Value2 = "M2"
$ws.Range("E23").Value2 = 3
$ws.Range("F23").Value2 = 1
$ws.Range("G23").Value2 = 6.849529333333333
$ws.Range("H23").Value2 = 20.548588
$ws.Range("I23").Value2 = 0.02599990892953593
$ws.Range("J23").Value2 = 0.03541161817534047
$ws.Range("K23").Value2 = 3
$ws.Range("L23").Value2 = 1
$ws.Range("M23").Value2 = 0.1136856666666667
$ws.Range("N23").Value2 = 0.341057
$ws.Range("O23").Value2 = 0.0004616453666923254
$ws.Range("P23").Value2 = 0.0005490449607626364
$ws.Range("Q23").Value2 = 0.7786933086128889
$ws.Range("R23").Value2 = 7.008239777516
$ws.Range("S23").Value2 = 0.00001200273749174268
$ws.Range("T23").Value2 = 0.00001944257051162127

# Row 24
$ws.Range("A24").Value2 = "M2"
$ws.Range("B24").Value2 = "Dag1"
$ws.Range("C24").Value2 = "Lama4"
$ws.Range("D24").Value2 = "Neutro"
$ws.Range("E24").Value2 = 3
$ws.Range("F24").Value2 = 1
$ws.Range("G24").Value2 = 6.849529333333333
$ws.Range("H24").Value2 = 20.548588
$ws.Range("I24").Value2 = 0.02599990892953593
$ws.Range("J24").Value2 = 0.03541161817534047
$ws.Range("K24").Value2 = 3
$ws.Range("L24").Value2 = 1
$ws.Range("M24").Value2 = 70.30838766666666
$ws.Range("N24").Value2 = 210.925163
$ws.Range("O24").Value2 = 0.285502494356584
$ws.Range("P24").Value2 = 0.3395543790134426
$ws.Range("Q24").Value2 = 481.5793637022048
$ws.Range("R24").Value2 = 4334.214273319843
$ws.Range("S24").Value2 = 0.00742303885242653
$ws.Range("T24").Value2 = 0.01202417001938887

# Row 25
$ws.Range("A25").Value2 = "M2"
$ws.Range("B25").Value2 = "Dag1"
$ws.Range("C25").Value2 = "Lama4"
$ws.Range("D25").Value2 = "sCs"
$ws.Range("E25").Value2 = 3
$ws.Range("F25").Value2 = 1
$ws.Range("G25").Value2 = 6.849529333333333
$ws.Range("H25").Value2 = 20.548588
$ws.Range("I25").Value2 = 0.02599990892953593
$ws.Range("J25").Value2 = 0.03541161817534047
$ws.Range("K25").Value2 = 2
$ws.Range("L25").Value2 = 1
$ws.Range("M25").Value2 = 13.5936885
$ws.Range("N25").Value2 = 27.187377
$ws.Range("O25").Value2 = 0.05520012765271271
$ws.Range("P25").Value2 = 0.04376714841860452
$ws.Range("Q25").Value2 = 93.11036812894599
$ws.Range("R25").Value2 = 558.662208773676
$ws.Range("S25").Value2 = 0.001435198291869288
$ws.Range("T25").Value2 = 0.00154986554842308

# Row 26
$ws.Range("A26").Value2 = "Neutro"
$ws.Range("B26").Value2 = "Dag1"
$ws.Range("C26").Value2 = "Lama4"
$ws.Range("D26").Value2 = "ECs"
$ws.Range("E26").Value2 = 3
$ws.Range("F26").Value2 = 1
$ws.Range("G26").Value2 = 15.48994666666667
$ws.Range("H26").Value2 = 46.46984
$ws.Range("I26").Value2 = 0.05879779223614324
$ws.Range("J26").Value2 = 0.08008201005096623
$ws.Range("K26").Value2 = 2
$ws.Range("L26").Value2 = 1
$ws.Range("M26").Value2 = 104.0097425
$ws.Range("N26").Value2 = 208.019485
$ws.Range("O26").Value2 = 0.4223541729035338
$ws.Range("P26").Value2 = 0.3348767214268842
$ws.Range("Q26").Value2 = 1611.105364138733
$ws.Range("R26").Value2 = 9666.6321848324
$ws.Range("S26").Value2 = 0.0248334929084501
$ws.Range("T26").Value2 = 0.02681760097114235

# Row 27
$ws.Range("A27").Value2 = "Neutro"
$ws.Range("B27").Value2 = "Dag1"
$ws.Range("C27").Value2 = "Lama4"
$ws.Range("D27").Value2 = "FAPs"
$ws.Range("E27").Value2 = 3
$ws.Range("F27").Value2 = 1
$ws.Range("G27").Value2 = 15.48994666666667
$ws.Range("H27").Value2 = 46.46984
$ws.Range("I27").Value2 = 0.05879779223614324
$ws.Range("J27").Value2 = 0.08008201005096623
$ws.Range("K27").Value2 = 3
$ws.Range("L27").Value2 = 1
$ws.Range("M27").Value2 = 58.20636966666666
$ws.Range("N27").Value2 = 174.619109
$ws.Range("O27").Value2 = 0.2363596190835905
$ws.Range("P27").Value2 = 0.2811076795060988
$ws.Range("Q27").Value2 = 901.6135617969511
$ws.Range("R27").Value2 = 8114.52205617256
$ws.Range("S27").Value2 = 0.01389742377589091
$ws.Range("T27").Value2 = 0.02251166801561119

# Row 28
$ws.Range("A28").Value2 = "Neutro"
$ws.Range("B28").Value2 = "Dag1"
$ws.Range("C28").Value2 = "Lama4"
$ws.Range("D28").Value2 = "M1"
$ws.Range("E28").Value2 = 3
$ws.Range("F28").Value2 = 1
$ws.Range("G28").Value2 = 15.48994666666667
$ws.Range("H28").Value2 = 46.46984
$ws.Range("I28").Value2 = 0.05879779223614324
$ws.Range("J28").Value2 = 0.08008201005096623
$ws.Range("K28").Value2 = 2
$ws.Range("L28").Value2 = 0.6666666666666666
$ws.Range("M28").Value2 = 0.03002933333333333
$ws.Range("N28").Value2 = 0.090088
$ws.Range("O28").Value2 = 0.0001219406368864389
$ws.Range("P28").Value2 = 0.0001450266742074914
$ws.Range("Q28").Value2 = 0.465152771768889
$ws.Range("R28").Value2 = 4.186374945920001
$ws.Range("S28").Value2 = 0.000007169840232791821
$ws.Range("T28").Value2 = 0.00001161402758154253

# Row 29
$ws.Range("A29").Value2 = "Neutro"
$ws.Range("B29").Value2 = "Dag1"
$ws.Range("C29").Value2 = "Lama4"
$ws.Range("D29").Value2 = "M2"
$ws.Range("E29").Value2 = 3
$ws.Range("F29").Value2 = 1
$ws.Range("G29").Value2 = 15.48994666666667
$ws.Range("H29").Value2 = 46.46984
$ws.Range("I29").Value2 = 0.05879779223614324
$ws.Range("J29").Value2 = 0.08008201005096623
$ws.Range("K29").Value2 = 3
$ws.Range("L29").Value2 = 1
$ws.Range("M29").Value2 = 0.1136856666666667
$ws.Range("N29").Value2 = 0.341057
$ws.Range("O29").Value2 = 0.0004616453666923254
$ws.Range("P29").Value2 = 0.0005490449607626364
$ws.Range("Q29").Value2 = 1.760984913431111
$ws.Range("R29").Value2 = 15.84886422088
$ws.Range("S29").Value2 = 0.00002714372835755351
$ws.Range("T29").Value2 = 0.0000439686240662258

# Row 30
$ws.Range("A30").Value2 = "Neutro"
$ws.Range("B30").Value2 = "Dag1"
$ws.Range("C30").Value2 = "Lama4"
$ws.Range("D30").Value2 = "Neutro"
$ws.Range("E30").Value2 = 3
$ws.Range("F30").Value2 = 1
$ws.Range("G30").Value2 = 15.48994666666667
$ws.Range("H30").Value2 = 46.46984
$ws.Range("I30").Value2 = 0.05879779223614324
$ws.Range("J30").Value2 = 0.08008201005096623
$ws.Range("K30").Value2 = 3
$ws.Range("L30").Value2 = 1
$ws.Range("M30").Value2 = 70.30838766666666
$ws.Range("N30").Value2 = 210.925163
$ws.Range("O30").Value2 = 0.285502494356584
$ws.Range("P30").Value2 = 0.3395543790134426
$ws.Range("Q30").Value2 = 1089.073175175991
$ws.Range("R30").Value2 = 9801.65857658392
$ws.Range("S30").Value2 = 0.01678691634607908
$ws.Range("T30").Value2 = 0.02719219719300411

# Row 31
$ws.Range("A31").Value2 = "Neutro"
$ws.Range("B31").Value2 = "Dag1"
$ws.Range("C31").Value2 = "Lama4"
$ws.Range("D31").Value2 = "sCs"
$ws.Range("E31").Value2 = 3
$ws.Range("F31").Value2 = 1
$ws.Range("G31").Value2 = 15.48994666666667
$ws.Range("H31").Value2 = 46.46984
$ws.Range("I31").Value2 = 0.05879779223614324
$ws.Range("J31").Value2 = 0.08008201005096623
$ws.Range("K31").Value2 = 2
$ws.Range("L31").Value2 = 1
$ws.Range("M31").Value2 = 13.5936885
$ws.Range("N31").Value2 = 27.187377
$ws.Range("O31").Value2 = 0.05520012765271271
$ws.Range("P31").Value2 = 0.04376714841860452
$ws.Range("Q31").Value2 = 210.56550986828
$ws.Range("R31").Value2 = 1263.39305920968
$ws.Range("S31").Value2 = 0.003245645637132787
$ws.Range("T31").Value2 = 0.003504961219560818

# Row 32
$ws.Range("A32").Value2 = "sCs"
$ws.Range("B32").Value2 = "Dag1"
$ws.Range("C32").Value2 = "Lama4"
$ws.Range("D32").Value2 = "ECs"
$ws.Range("E32").Value2 = 2
$ws.Range("F32").Value2 = 1
$ws.Range("G32").Value2 = 189.0471995
$ws.Range("H32").Value2 = 378.094399
$ws.Range("I32").Value2 = 0.717598207290517
$ws.Range("J32").Value2 = 0.6515744289399755
$ws.Range("K32").Value2 = 2
$ws.Range("L32").Value2 = 1
$ws.Range("M32").Value2 = 104.0097425
$ws.Range("N32").Value2 = 208.019485
$ws.Range("O32").Value2 = 0.4223541729035338
$ws.Range("P32").Value2 = 0.3348767214268842
$ws.Range("Q32").Value2 = 19662.75054034113
$ws.Range("R32").Value2 = 78651.00216136452
$ws.Range("S32").Value2 = 0.3030805973172449
$ws.Range("T32").Value2 = 0.2181971085290133

# Row 33
$ws.Range("A33").Value2 = "sCs"
$ws.Range("B33").Value2 = "Dag1"
$ws.Range("C33").Value2 = "Lama4"
$ws.Range("D33").Value2 = "FAPs"
$ws.Range("E33").Value2 = 2
$ws.Range("F33").Value2 = 1
$ws.Range("G33").Value2 = 189.0471995
$ws.Range("H33").Value2 = 378.094399
$ws.Range("I33").Value2 = 0.717598207290517
$ws.Range("J33").Value2 = 0.6515744289399755
$ws.Range("K33").Value2 = 3
$ws.Range("L33").Value2 = 1
$ws.Range("M33").Value2 = 58.20636966666666
$ws.Range("N33").Value2 = 174.619109
$ws.Range("O33").Value2 = 0.2363596190835905
$ws.Range("P33").Value2 = 0.2811076795060988
$ws.Range("Q33").Value2 = 11003.75117854508
$ws.Range("R33").Value2 = 66022.50707127049
$ws.Range("S33").Value2 = 0.169611238930254
$ws.Range("T33").Value2 = 0.183162575744828

# Row 34
$ws.Range("A34").Value2 = "sCs"
$ws.Range("B34").Value2 = "Dag1"
$ws.Range("C34").Value2 = "Lama4"
$ws.Range("D34").Value2 = "M1"
$ws.Range("E34").Value2 = 2
$ws.Range("F34").Value2 = 1
$ws.Range("G34").Value2 = 189.0471995
$ws.Range("H34").Value2 = 378.094399
$ws.Range("I34").Value2 = 0.717598207290517
$ws.Range("J34").Value2 = 0.6515744289399755
$ws.Range("K34").Value2 = 2
$ws.Range("L34").Value2 = 0.6666666666666666
$ws.Range("M34").Value2 = 0.03002933333333333
$ws.Range("N34").Value2 = 0.090088
$ws.Range("O34").Value2 = 0.0001219406368864389
$ws.Range("P34").Value2 = 0.0001450266742074914
$ws.Range("Q34").Value2 = 5.676961369518667
$ws.Range("R34").Value2 = 34.061768217112
$ws.Range("S34").Value2 = 0.00008750438242557249
$ws.Range("T34").Value2 = 0.00009449567242781007

# Row 35
$ws.Range("A35").Value2 = "sCs"
$ws.Range("B35").Value2 = "Dag1"
$ws.Range("C35").Value2 = "Lama4"
$ws.Range("D35").Value2 = "M2"
$ws.Range("E35").Value2 = 2
$ws.Range("F35").Value2 = 1
$ws.Range("G35").Value2 = 189.0471995
$ws.Range("H35").Value2 = 378.094399
$ws.Range("I35").Value2 = 0.717598207290517
$ws.Range("J35").Value2 = 0.6515744289399755
$ws.Range("K35").Value2 = 3
$ws.Range("L35").Value2 = 1
$ws.Range("M35").Value2 = 0.1136856666666667
$ws.Range("N35").Value2 = 0.341057
$ws.Range("O35").Value2 = 0.0004616453666923254
$ws.Range("P35").Value2 = 0.0005490449607626364
$ws.Range("Q35").Value2 = 21.49195690662383
$ws.Range("R35").Value2 = 128.951741439743
$ws.Range("S35").Value2 = 0.0003312758875423861
$ws.Range("T35").Value2 = 0.0003577436567712861

# Row 36
$ws.Range("A36").Value2 = "sCs"
$ws.Range("B36").Value2 = "Dag1"
$ws.Range("C36").Value2 = "Lama4"
$ws.Range("D36").Value2 = "Neutro"
$ws.Range("E36").Value2 = 2
$ws.Range("F36").Value2 = 1
$ws.Range("G36").Value2 = 189.0471995
$ws.Range("H36").Value2 = 378.094399
$ws.Range("I36").Value2 = 0.717598207290517
$ws.Range("J36").Value2 = 0.6515744289399755
$ws.Range("K36").Value2 = 3
$ws.Range("L36").Value2 = 1
$ws.Range("M36").Value2 = 70.30838766666666
$ws.Range("N36").Value2 = 210.925163
$ws.Range("O36").Value2 = 0.285502494356584
$ws.Range("P36").Value2 = 0.3395543790134426
$ws.Range("Q36").Value2 = 13291.60378974367
$ws.Range("R36").Value2 = 79749.62273846204
$ws.Range("S36").Value2 = 0.2048760781272557
$ws.Range("T36").Value2 = 0.2212449505997519

# Row 37
$ws.Range("A37").Value2 = "sCs"
$ws.Range("B37").Value2 = "Dag1"
$ws.Range("C37").Value2 = "Lama4"
$ws.Range("D37").Value2 = "sCs"
$ws.Range("E37").Value2 = 2
$ws.Range("F37").Value2 = 1
$ws.Range("G37").Value2 = 189.0471995
$ws.Range("H37").Value2 = 378.094399
$ws.Range("I37").Value2 = 0.717598207290517
$ws.Range("J37").Value2 = 0.6515744289399755
$ws.Range("K37").Value2 = 2
$ws.Range("L37").Value2 = 1
$ws.Range("M37").Value2 = 13.5936885
$ws.Range("N37").Value2 = 27.187377
$ws.Range("O37").Value2 = 0.05520012765271271
$ws.Range("P37").Value2 = 0.04376714841860452
$ws.Range("Q37").Value2 = 2569.848741800356
$ws.Range("R37").Value2 = 10279.39496720142
$ws.Range("S37").Value2 = 0.03961151264579434
$ws.Range("T37").Value2 = 0.02851755473718339
